# Update the Hicks-Henne deformed-airfoil coordinate table (Sheet1, A:B)
# with the refreshed (finer, 81-point) x/y_modified dataset and extend
# the used range from A1:B36 to A1:B82.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 0.00126
$ws.Cells.Item(3, 1).Value = 0.998459
$ws.Cells.Item(3, 2).Value = 0.001476333889201003
$ws.Cells.Item(4, 1).Value = 0.9938439999999999
$ws.Cells.Item(4, 2).Value = 0.002125337068260748
$ws.Cells.Item(5, 1).Value = 0.986185
$ws.Cells.Item(5, 2).Value = 0.003208948377701005
$ws.Cells.Item(6, 1).Value = 0.975528
$ws.Cells.Item(6, 2).Value = 0.004726845787834165
$ws.Cells.Item(7, 1).Value = 0.96194
$ws.Cells.Item(7, 2).Value = 0.00668402119748149
$ws.Cells.Item(8, 1).Value = 0.945503
$ws.Cells.Item(8, 2).Value = 0.00908211523504465
$ws.Cells.Item(9, 1).Value = 0.92632
$ws.Cells.Item(9, 2).Value = 0.01192724388427125
$ws.Cells.Item(10, 1).Value = 0.904508
$ws.Cells.Item(10, 2).Value = 0.01522544591879739
$ws.Cells.Item(11, 1).Value = 0.880203
$ws.Cells.Item(11, 2).Value = 0.01898232102286982
$ws.Cells.Item(12, 1).Value = 0.853553
$ws.Cells.Item(12, 2).Value = 0.02319926111471215
$ws.Cells.Item(13, 1).Value = 0.824724
$ws.Cells.Item(13, 2).Value = 0.02787363714074887
$ws.Cells.Item(14, 1).Value = 0.793893
$ws.Cells.Item(14, 2).Value = 0.03299086005392982
$ws.Cells.Item(15, 1).Value = 0.761249
$ws.Cells.Item(15, 2).Value = 0.03852384859357359
$ws.Cells.Item(16, 1).Value = 0.7269949999999999
$ws.Cells.Item(16, 2).Value = 0.04442630530896738
$ws.Cells.Item(17, 1).Value = 0.691342
$ws.Cells.Item(17, 2).Value = 0.05062789965516801
$ws.Cells.Item(18, 1).Value = 0.654508
$ws.Cells.Item(18, 2).Value = 0.05703111215392148
$ws.Cells.Item(19, 1).Value = 0.616723
$ws.Cells.Item(19, 2).Value = 0.06351093398096402
$ws.Cells.Item(20, 1).Value = 0.578217
$ws.Cells.Item(20, 2).Value = 0.06991083669794433
$ws.Cells.Item(21, 1).Value = 0.53923
$ws.Cells.Item(21, 2).Value = 0.07604921899626271
$ws.Cells.Item(22, 1).Value = 0.5
$ws.Cells.Item(22, 2).Value = 0.08172417134448536
$ws.Cells.Item(23, 1).Value = 0.46077
$ws.Cells.Item(23, 2).Value = 0.08671918650178619
$ws.Cells.Item(24, 1).Value = 0.421783
$ws.Cells.Item(24, 2).Value = 0.09081421746202267
$ws.Cells.Item(25, 1).Value = 0.383277
$ws.Cells.Item(25, 2).Value = 0.09380337452555637
$ws.Cells.Item(26, 1).Value = 0.345492
$ws.Cells.Item(26, 2).Value = 0.09550514476739626
$ws.Cells.Item(27, 1).Value = 0.308658
$ws.Cells.Item(27, 2).Value = 0.09577550572467908
$ws.Cells.Item(28, 1).Value = 0.273005
$ws.Cells.Item(28, 2).Value = 0.09452645363743438
$ws.Cells.Item(29, 1).Value = 0.238751
$ws.Cells.Item(29, 2).Value = 0.09173446030778616
$ws.Cells.Item(30, 1).Value = 0.206107
$ws.Cells.Item(30, 2).Value = 0.0874464056850444
$ws.Cells.Item(31, 1).Value = 0.175276
$ws.Cells.Item(31, 2).Value = 0.08178522610095247
$ws.Cells.Item(32, 1).Value = 0.146447
$ws.Cells.Item(32, 2).Value = 0.07494342737709588
$ws.Cells.Item(33, 1).Value = 0.119797
$ws.Cells.Item(33, 2).Value = 0.06717454024458965
$ws.Cells.Item(34, 1).Value = 0.09549199999999999
$ws.Cells.Item(34, 2).Value = 0.05877525087092676
$ws.Cells.Item(35, 1).Value = 0.07368
$ws.Cells.Item(35, 2).Value = 0.05006512483005968
$ws.Cells.Item(36, 1).Value = 0.054497
$ws.Cells.Item(36, 2).Value = 0.04136231091921833
$ws.Cells.Item(37, 1).Value = 0.03806
$ws.Cells.Item(37, 2).Value = 0.03294978193497399
$ws.Cells.Item(38, 1).Value = 0.024472
$ws.Cells.Item(38, 2).Value = 0.02505861721212721
$ws.Cells.Item(39, 1).Value = 0.013815
$ws.Cells.Item(39, 2).Value = 0.01784941770713892
$ws.Cells.Item(40, 1).Value = 0.006156
$ws.Cells.Item(40, 2).Value = 0.01141753350478415
$ws.Cells.Item(41, 1).Value = 0.001541
$ws.Cells.Item(41, 2).Value = 0.005872683264668331
$ws.Cells.Item(42, 1).Value = 0
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(43, 1).Value = 0.001541
$ws.Cells.Item(43, 2).Value = -0.007916731751757536
$ws.Cells.Item(44, 1).Value = 0.006156
$ws.Cells.Item(44, 2).Value = -0.01627153601609284
$ws.Cells.Item(45, 1).Value = 0.013815
$ws.Cells.Item(45, 2).Value = -0.02412849709361464
$ws.Cells.Item(46, 1).Value = 0.024472
$ws.Cells.Item(46, 2).Value = -0.03165366158023353
$ws.Cells.Item(47, 1).Value = 0.03806
$ws.Cells.Item(47, 2).Value = -0.03886242317256935
$ws.Cells.Item(48, 1).Value = 0.054497
$ws.Cells.Item(48, 2).Value = -0.04570804941243417
$ws.Cells.Item(49, 1).Value = 0.07368
$ws.Cells.Item(49, 2).Value = -0.05211557295782788
$ws.Cells.Item(50, 1).Value = 0.09549199999999999
$ws.Cells.Item(50, 2).Value = -0.0579964732709505
$ws.Cells.Item(51, 1).Value = 0.119797
$ws.Cells.Item(51, 2).Value = -0.06324902010494303
$ws.Cells.Item(52, 1).Value = 0.146447
$ws.Cells.Item(52, 2).Value = -0.06777191408801952
$ws.Cells.Item(53, 1).Value = 0.175276
$ws.Cells.Item(53, 2).Value = -0.07147244367126217
$ws.Cells.Item(54, 1).Value = 0.206107
$ws.Cells.Item(54, 2).Value = -0.07427315300689262
$ws.Cells.Item(55, 1).Value = 0.238751
$ws.Cells.Item(55, 2).Value = -0.07612169983841448
$ws.Cells.Item(56, 1).Value = 0.273005
$ws.Cells.Item(56, 2).Value = -0.07699540246202502
$ws.Cells.Item(57, 1).Value = 0.308658
$ws.Cells.Item(57, 2).Value = -0.07690659823547313
$ws.Cells.Item(58, 1).Value = 0.345492
$ws.Cells.Item(58, 2).Value = -0.07589912858511765
$ws.Cells.Item(59, 1).Value = 0.383277
$ws.Cells.Item(59, 2).Value = -0.0740464183016423
$ws.Cells.Item(60, 1).Value = 0.421783
$ws.Cells.Item(60, 2).Value = -0.07144811014533453
$ws.Cells.Item(61, 1).Value = 0.46077
$ws.Cells.Item(61, 2).Value = -0.0682177379351675
$ws.Cells.Item(62, 1).Value = 0.5
$ws.Cells.Item(62, 2).Value = -0.06447630728366603
$ws.Cells.Item(63, 1).Value = 0.53923
$ws.Cells.Item(63, 2).Value = -0.06034871150286349
$ws.Cells.Item(64, 1).Value = 0.578217
$ws.Cells.Item(64, 2).Value = -0.05595269367257467
$ws.Cells.Item(65, 1).Value = 0.616723
$ws.Cells.Item(65, 2).Value = -0.05139412998728005
$ws.Cells.Item(66, 1).Value = 0.654508
$ws.Cells.Item(66, 2).Value = -0.04676628031488771
$ws.Cells.Item(67, 1).Value = 0.691342
$ws.Cells.Item(67, 2).Value = -0.04214978350412942
$ws.Cells.Item(68, 1).Value = 0.7269949999999999
$ws.Cells.Item(68, 2).Value = -0.03760861837390687
$ws.Cells.Item(69, 1).Value = 0.761249
$ws.Cells.Item(69, 2).Value = -0.03319618697557988
$ws.Cells.Item(70, 1).Value = 0.793893
$ws.Cells.Item(70, 2).Value = -0.02895544663309145
$ws.Cells.Item(71, 1).Value = 0.824724
$ws.Cells.Item(71, 2).Value = -0.02492098726715349
$ws.Cells.Item(72, 1).Value = 0.853553
$ws.Cells.Item(72, 2).Value = -0.02112195903756597
$ws.Cells.Item(73, 1).Value = 0.880203
$ws.Cells.Item(73, 2).Value = -0.01758594305949756
$ws.Cells.Item(74, 1).Value = 0.904508
$ws.Cells.Item(74, 2).Value = -0.01433644197668487
$ws.Cells.Item(75, 1).Value = 0.92632
$ws.Cells.Item(75, 2).Value = -0.01139796148044112
$ws.Cells.Item(76, 1).Value = 0.945503
$ws.Cells.Item(76, 2).Value = -0.008792934797641662
$ws.Cells.Item(77, 1).Value = 0.96194
$ws.Cells.Item(77, 2).Value = -0.006543270889450678
$ws.Cells.Item(78, 1).Value = 0.975528
$ws.Cells.Item(78, 2).Value = -0.004668795992718405
$ws.Cells.Item(79, 1).Value = 0.986185
$ws.Cells.Item(79, 2).Value = -0.003190491699260172
$ws.Cells.Item(80, 1).Value = 0.9938439999999999
$ws.Cells.Item(80, 2).Value = -0.002121679248347373
$ws.Cells.Item(81, 1).Value = 0.998459
$ws.Cells.Item(81, 2).Value = -0.001476104965016639
$ws.Cells.Item(82, 1).Value = 1
$ws.Cells.Item(82, 2).Value = -0.00126
